$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the full Area / Subarea names with their short codes
$ws.Range("E2").Value = "BC"
$ws.Range("F2").Value = "PAR"
$ws.Range("H2").Value = "BQM"

$ws.Range("E3").Value = "BC"
$ws.Range("F3").Value = "PAR"
$ws.Range("H3").Value = "BQM"

# Restore view to top-left and select F3 (matches author's final cursor position)
$ws.Range("A1").Select()
$ws.Range("F3").Select()
